$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8836637139320374
$ws.Range("B1").Value = 1.036612510681152
$ws.Range("C1").Value = 0.8401211500167847
$ws.Range("D1").Value = 0.7816460132598877
$ws.Range("E1").Value = 0.8130704760551453
